$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on the phone-number cells so the numeric-looking
# values keep their leading zeros instead of being auto-converted to numbers
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"

# Row 2 - update guest entry (name stays the same, other fields become text values)
$ws.Range("C2").Value = "Mega Siva Marhaeni"
$ws.Range("D2").Value = "Kampung Babakan Sate"
$ws.Range("E2").Value = "08080808"
$ws.Range("F2").Value = "Nur Syifa Assawwala"
$ws.Range("G2").Value = "Nyeblak di Teh Anna"

# Row 3 - update guest entry
$ws.Range("C3").Value = "Muhammad Haidar Almer Rafif"
$ws.Range("D3").Value = "Perumahan Bumi Marhamah"
$ws.Range("E3").Value = "085956267079"
$ws.Range("F3").Value = "Presiden"
$ws.Range("G3").Value = "Reformasi"
